$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column after each of the 10 original 2-column stock blocks
# (processed right-to-left so column letters stay valid while inserting)
$insertBeforeCols = @("U","S","Q","O","M","K","I","G","E","C")
foreach ($col in $insertBeforeCols) {
    $ws.Range($col + "1").EntireColumn.Insert()
}

# Append two more new columns at the end: one finishes the 11th block (CHINHIN),
# the other is the entirely new 12th stock (KRETAM)
$ws.Range("AG1").EntireColumn.Insert()
$ws.Range("AH1").EntireColumn.Insert()

# Fill in the newly inserted columns (2020-03-12 data, plus the new KRETAM column in AH)
$ws.Range("C1").Value = 2
$ws.Range("F1").Value = 2
$ws.Range("I1").Value = 2
$ws.Range("L1").Value = 2
$ws.Range("O1").Value = 2
$ws.Range("R1").Value = 2
$ws.Range("U1").Value = 2
$ws.Range("X1").Value = 2
$ws.Range("AA1").Value = 2
$ws.Range("AD1").Value = 2
$ws.Range("AG1").Value = 2
$ws.Range("AH1").Value = 0

$ws.Range("C2").Value = 1583971200
$ws.Range("F2").Value = 1583971200
$ws.Range("I2").Value = 1583971200
$ws.Range("L2").Value = 1583971200
$ws.Range("O2").Value = 1583971200
$ws.Range("R2").Value = 1583971200
$ws.Range("U2").Value = 1583971200
$ws.Range("X2").Value = 1583971200
$ws.Range("AA2").Value = 1583971200
$ws.Range("AD2").Value = 1583971200
$ws.Range("AG2").Value = 1583971200
$ws.Range("AH2").Value = 1583971200

$ws.Range("C3").Value = "'2020-03-12"
$ws.Range("F3").Value = "'2020-03-12"
$ws.Range("I3").Value = "'2020-03-12"
$ws.Range("L3").Value = "'2020-03-12"
$ws.Range("O3").Value = "'2020-03-12"
$ws.Range("R3").Value = "'2020-03-12"
$ws.Range("U3").Value = "'2020-03-12"
$ws.Range("X3").Value = "'2020-03-12"
$ws.Range("AA3").Value = "'2020-03-12"
$ws.Range("AD3").Value = "'2020-03-12"
$ws.Range("AG3").Value = "'2020-03-12"
$ws.Range("AH3").Value = "'2020-03-12"

$ws.Range("C4").Value = "'4723"
$ws.Range("F4").Value = "'0083"
$ws.Range("I4").Value = "'9466"
$ws.Range("L4").Value = "'0215"
$ws.Range("O4").Value = "'5277"
$ws.Range("R4").Value = "'5292"
$ws.Range("U4").Value = "'0208"
$ws.Range("X4").Value = "'0176"
$ws.Range("AA4").Value = "'0198"
$ws.Range("AD4").Value = "'0128"
$ws.Range("AG4").Value = "'5273"
$ws.Range("AH4").Value = "'1996"

$ws.Range("C5").Value = "'JAKS"
$ws.Range("F5").Value = "'NOTION"
$ws.Range("I5").Value = "'KKB"
$ws.Range("L5").Value = "'SLVEST"
$ws.Range("O5").Value = "'FPGROUP"
$ws.Range("R5").Value = "'UWC"
$ws.Range("U5").Value = "'GREATEC"
$ws.Range("X5").Value = "'KRONO"
$ws.Range("AA5").Value = "'GDB"
$ws.Range("AD5").Value = "'FRONTKN"
$ws.Range("AG5").Value = "'CHINHIN"
$ws.Range("AH5").Value = "'KRETAM"

$ws.Range("C6").Value = 0.985
$ws.Range("F6").Value = 0.775
$ws.Range("I6").Value = 1.86
$ws.Range("L6").Value = 0.875
$ws.Range("O6").Value = 0.67
$ws.Range("R6").Value = 2.08
$ws.Range("U6").Value = 2.8
$ws.Range("X6").Value = 0.5
$ws.Range("AA6").Value = 0.66
$ws.Range("AD6").Value = 1.91
$ws.Range("AG6").Value = 0.5600000000000001
$ws.Range("AH6").Value = 0.315

$ws.Range("C7").Value = 0.99
$ws.Range("F7").Value = 0.78
$ws.Range("I7").Value = 1.86
$ws.Range("L7").Value = 0.885
$ws.Range("O7").Value = 0.675
$ws.Range("R7").Value = 2.09
$ws.Range("U7").Value = 2.8
$ws.Range("X7").Value = 0.505
$ws.Range("AA7").Value = 0.66
$ws.Range("AD7").Value = 1.95
$ws.Range("AG7").Value = 0.575
$ws.Range("AH7").Value = 0.38

$ws.Range("C8").Value = 0.865
$ws.Range("F8").Value = 0.735
$ws.Range("I8").Value = 1.76
$ws.Range("L8").Value = 0.8
$ws.Range("O8").Value = 0.615
$ws.Range("R8").Value = 1.83
$ws.Range("U8").Value = 2.47
$ws.Range("X8").Value = 0.47
$ws.Range("AA8").Value = 0.63
$ws.Range("AD8").Value = 1.86
$ws.Range("AG8").Value = 0.55
$ws.Range("AH8").Value = 0.3

$ws.Range("C9").Value = 0.87
$ws.Range("F9").Value = 0.74
$ws.Range("I9").Value = 1.77
$ws.Range("L9").Value = 0.8100000000000001
$ws.Range("O9").Value = 0.625
$ws.Range("R9").Value = 1.87
$ws.Range("U9").Value = 2.47
$ws.Range("X9").Value = 0.47
$ws.Range("AA9").Value = 0.64
$ws.Range("AD9").Value = 1.89
$ws.Range("AG9").Value = 0.575
$ws.Range("AH9").Value = 0.375

$ws.Range("C10").Value = 882187
$ws.Range("F10").Value = 79996
$ws.Range("I10").Value = 19681
$ws.Range("L10").Value = 131196
$ws.Range("O10").Value = 212709
$ws.Range("R10").Value = 229494
$ws.Range("U10").Value = 86150
$ws.Range("X10").Value = 103795
$ws.Range("AA10").Value = 48231
$ws.Range("AD10").Value = 92928
$ws.Range("AG10").Value = 50667
$ws.Range("AH10").Value = 33240

